$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 22 ("Update author" use case), pushing the
# "Get all publisher" / "Create new publisher" rows (and everything below)
# down by one row. Excel copies the formatting of the row above (row 21)
# into the freshly inserted row.
$ws.Rows("22:22").Insert()

# Fill in the new use-case row.
$ws.Range("C22").Value = "Update"
$ws.Range("D22").Value = "Yes"
$ws.Range("E22").Value = "[POST]book/update_author"
$ws.Range("F22").Value = "{`n     ""id"": 7,`n     ""name"": ""Nam Cao"",`n     ""currentUserID"": 1`n}"
$ws.Range("G22").Value = "{`n    ""result"": true,`n    ""message"": ""Cập nhật thành công"",`n    ""data"": {`n        ""id"": 7,`n        ""name"": ""Nam Cao"",`n        ""updatedDate"": ""15-08-2020"",`n        ""updatedAccount"": ""Võ Thanh Hiếu"",`n        ""updatedAccountID"": 1`n    }`n}"

# A22/B22 stay empty (no STT / no use-case group label), but centred both
# ways like the rest of the numbered rows.
$rng = $ws.Range("A22:B22")
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

# The row holds a lot of wrapped text, same row height as the other
# JSON request/response rows.
$ws.Rows("22:22").RowHeight = 165

# Update the saved view/selection to where the author was last looking.
$ws.Range("F16").Select()
$ws.Range("H24").Select()
